$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A45").Value = "How many tables can I have in my log?"
$ws.Range("B45").Value = "You can have up to 100 tables in a log."
